$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 128, shifting the existing rows 128:146 down to 129:147.
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with the new outbound-delivery record.
$ws.Cells.Item(128, 1).Value = "50170625"
$ws.Cells.Item(128, 2).Value = "10"
$ws.Cells.Item(128, 3).Value = "07019-0"
$ws.Cells.Item(128, 4).Value = "CONTACTOR 500VDC 535AMP 5DP8-5021-21"
$ws.Cells.Item(128, 5).Value = 1
$ws.Cells.Item(128, 6).Value = "PC"
$ws.Cells.Item(128, 7).Value = "10122213"
$ws.Cells.Item(128, 8).Value = "1010027"
$ws.Cells.Item(128, 9).Value = "Completed"
$ws.Cells.Item(128, 10).Value = [DateTime]"2026-01-21"
$ws.Cells.Item(128, 11).Value = [DateTime]"2026-01-28"
$ws.Cells.Item(128, 12).Value = [DateTime]"2026-01-28"
$ws.Cells.Item(128, 13).Value = "ODLV"
$ws.Cells.Item(128, 14).Value = "Standard Item - Outbound Delivery"
